$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 125, shifting existing rows 125..180 down to 126..181.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new data record.
$ws.Range("A125").Value = 4
$ws.Range("B125").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C125").Value = "Los Lagos"
$ws.Range("D125").Value = 44609
$ws.Range("E125").Value = 10
$ws.Range("F125").Value = 100112039
$ws.Range("G125").Value = "Ciboulette"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 80
$ws.Range("K125").Value = 2500
$ws.Range("L125").Value = 2500
$ws.Range("M125").Value = 2500
$ws.Range("N125").Value = "$/docena de atados"
$ws.Range("O125").Value = "Región Metropolitana"
$ws.Range("P125").Value = 833
$ws.Range("Q125").Value = 3
$ws.Range("R125").Value = "Hortaliza"
